$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the row-5 measurement values to 2 decimal places (custom accuracy).
$row5 = @(
    6.62, 4.82, 0.77, 14.62, 11.47, 5.15, 22.93, 8.14, 3.51, 5.01, 5.85,
    6.3, 1.69, 5.26, 7.4, 4.65, 0.68, 0.43, 72.8, 14.83, 4.86, 9.75,
    5.07, 0.99, 10.86, 4.29, 3.91, 4.59, 6.12, 0.54, 21.1, 2.63, 6.07
)

for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5[$i]
}

# Row 6 was dropped from the dataset; remove it entirely (shifts dimension
# from A1:AH6 down to A1:AH5 automatically).
$ws.Rows(6).Delete()
